$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 0

$ws.Range("F8").Value = 50

$ws.Range("F10").Value = 60

$ws.Range("F14").Value = 60

$ws.Range("F15").Value = 0

$ws.Range("D16").Value = 0

$ws.Range("D17").Value = 0

# Update the selection to match the final state (A19:G19, active cell A19)
$ws.Range("A19:G19").Select()
